# Updates TPM-derived LR-pair statistics (new TPM values) and
# corrects the "Resolving-Mac" / "MuSCs" sending-cluster label for rows 7-11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @("I2", 0.9878740142699144),
  @("J2", 0.9878740142699145),
  @("M2", 119.0164006666667),
  @("N2", 357.049202),
  @("O2", 0.9176278005170622),
  @("P2", 0.9176278005170622),
  @("Q2", 150.316642894394),
  @("R2", 1352.849786049546),
  @("S2", 0.9065006589024625),
  @("T2", 0.9065006589024626),
  @("I3", 0.9878740142699144),
  @("J3", 0.9878740142699145),
  @("O3", 0.001755838010330732),
  @("P3", 0.001755838010330731),
  @("S3", 0.001734546743673119),
  @("T3", 0.001734546743673119),
  @("I4", 0.9878740142699144),
  @("J4", 0.9878740142699145),
  @("M4", 7.816301333333333),
  @("N4", 23.448904),
  @("O4", 0.06026442877207647),
  @("P4", 0.06026442877207646),
  @("Q4", 9.871918237287998),
  @("R4", 88.847264135592),
  @("S4", 0.0595336631687545),
  @("T4", 0.0595336631687545),
  @("I5", 0.9878740142699144),
  @("J5", 0.9878740142699145),
  @("M5", 0.105045),
  @("N5", 0.315135),
  @("O5", 0.000809906968832672),
  @("P5", 0.000809906968832672),
  @("Q5", 0.132670889595),
  @("R5", 1.194038006355),
  @("S5", 0.00080008604848591),
  @("T5", 0.0008000860484859101),
  @("I6", 0.9878740142699144),
  @("J6", 0.9878740142699145),
  @("M6", 2.534602333333333),
  @("N6", 7.603807),
  @("O6", 0.01954202573169801),
  @("P6", 0.01954202573169801),
  @("Q6", 3.201179935579),
  @("R6", 28.810619420211),
  @("S6", 0.01930505940653848),
  @("T6", 0.01930505940653848),
  @("A7", "MuSCs"),
  @("G7", 0.015503),
  @("H7", 0.046509),
  @("I7", 0.01212598573008556),
  @("J7", 0.01212598573008556),
  @("M7", 119.0164006666667),
  @("N7", 357.049202),
  @("O7", 0.9176278005170622),
  @("P7", 0.9176278005170622),
  @("Q7", 1.845111259535333),
  @("R7", 16.606001335818),
  @("S7", 0.01112714161459969),
  @("T7", 0.01112714161459969),
  @("A8", "MuSCs"),
  @("G8", 0.015503),
  @("H8", 0.046509),
  @("I8", 0.01212598573008556),
  @("J8", 0.01212598573008556),
  @("O8", 0.001755838010330732),
  @("P8", 0.001755838010330731),
  @("Q8", 0.003530534363666667),
  @("R8", 0.03177480927300001),
  @("S8", 0.00002129126665761226),
  @("T8", 0.00002129126665761226),
  @("A9", "MuSCs"),
  @("G9", 0.015503),
  @("H9", 0.046509),
  @("I9", 0.01212598573008556),
  @("J9", 0.01212598573008556),
  @("M9", 7.816301333333333),
  @("N9", 23.448904),
  @("O9", 0.06026442877207647),
  @("P9", 0.06026442877207646),
  @("Q9", 0.1211761195706667),
  @("R9", 1.090585076136),
  @("S9", 0.0007307656033219566),
  @("T9", 0.0007307656033219565),
  @("A10", "MuSCs"),
  @("G10", 0.015503),
  @("H10", 0.046509),
  @("I10", 0.01212598573008556),
  @("J10", 0.01212598573008556),
  @("M10", 0.105045),
  @("N10", 0.315135),
  @("O10", 0.000809906968832672),
  @("P10", 0.000809906968832672),
  @("Q10", 0.001628512635),
  @("R10", 0.014656613715),
  @("S10", 0.000009820920346761827),
  @("T10", 0.000009820920346761827),
  @("A11", "MuSCs"),
  @("G11", 0.015503),
  @("H11", 0.046509),
  @("I11", 0.01212598573008556),
  @("J11", 0.01212598573008556),
  @("M11", 2.534602333333333),
  @("N11", 7.603807),
  @("O11", 0.01954202573169801),
  @("P11", 0.01954202573169801),
  @("Q11", 0.03929393997366667),
  @("R11", 0.353645459763),
  @("S11", 0.0002369663251595348),
  @("T11", 0.0002369663251595348)
)

foreach ($u in $updates) {
  $ws.Range($u[0]).Value = $u[1]
}
